# Updated symbol list on Sat Jan 28 17:26:42 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row with
# the latest scraped values. Cells are plain text (e.g. "306.69", "-0.15%"),
# so force a Text number format before assigning the literal string --
# otherwise Excel's COM layer would auto-coerce these into numeric/percent
# values instead of preserving the original text representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue "D2" "306.69"
Set-TextValue "E2" "-0.15%"
Set-TextValue "D3" "38.96"
Set-TextValue "E3" "7.35%"
Set-TextValue "D4" "5.105"
Set-TextValue "E4" "0.97%"
Set-TextValue "D5" "0.08053"
Set-TextValue "E5" "-0.30%"
Set-TextValue "D6" "1.925"
Set-TextValue "E6" "-4.59%"
Set-TextValue "D7" "4.197"
Set-TextValue "E7" "1.29%"
Set-TextValue "D8" "7.975"
Set-TextValue "E8" "1.77%"
Set-TextValue "D9" "0.9339"
Set-TextValue "E9" "0.92%"
Set-TextValue "D10" "0.1452"
Set-TextValue "E10" "-0.73%"
Set-TextValue "D11" "0.1925"
Set-TextValue "E11" "0.13%"
Set-TextValue "D12" "0.09104"
Set-TextValue "E12" "0.27%"
Set-TextValue "E13" "2.10%"
Set-TextValue "D14" "0.09792"
Set-TextValue "E14" "-1.26%"
Set-TextValue "D15" "0.001401"
Set-TextValue "E15" "-0.57%"
Set-TextValue "D16" "0.005924"
Set-TextValue "E16" "-4.49%"
Set-TextValue "D17" "3.795"
Set-TextValue "E17" "-1.17%"
Set-TextValue "D18" "3.459"
Set-TextValue "E18" "1.98%"
Set-TextValue "E19" "-0.18%"
Set-TextValue "D20" "0.1304"
Set-TextValue "E20" "-2.35%"
Set-TextValue "D21" "4.780"
Set-TextValue "E21" "-0.41%"
Set-TextValue "D22" "0.2510"
Set-TextValue "E22" "7.37%"
Set-TextValue "D23" "0.04392"
Set-TextValue "E23" "0.75%"
Set-TextValue "D24" "0.001239"
Set-TextValue "E24" "0.84%"
Set-TextValue "D25" "0.004270"
Set-TextValue "E25" "-0.75%"
Set-TextValue "D26" "0.0001301"
Set-TextValue "E26" "0.21%"
Set-TextValue "D39" "0.02042"
Set-TextValue "E39" "1.16%"
Set-TextValue "D40" "0.05038"
Set-TextValue "E40" "-2.11%"
Set-TextValue "D41" "0.007435"
Set-TextValue "E41" "-0.68%"
Set-TextValue "D42" "0.01013"
Set-TextValue "E42" "0.09%"
Set-TextValue "D43" "0.1348"
Set-TextValue "E43" "-1.09%"
Set-TextValue "D44" "0.002142"
Set-TextValue "E44" "-0.72%"
Set-TextValue "D45" "0.009076"
Set-TextValue "E45" "-8.57%"
Set-TextValue "D46" "0.00006201"
Set-TextValue "E46" "-1.12%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.35%"
Set-TextValue "D48" "0.002806"
Set-TextValue "E49" "28.26%"
Set-TextValue "D50" "0.00002105"
Set-TextValue "E50" "0.35%"
Set-TextValue "D51" "0.0002005"
Set-TextValue "E51" "0.35%"
